$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 14, column A text (Oct 2025 update: course now spans a full week) ---
# Free up the old shared string slot by first writing the new "genomic medicine" course name
# into B15 (this reuses the slot vacated once A14's old text is no longer referenced elsewhere),
# matching the order of edits the original author made.
$ws.Range("B15").Value = "MSc in Genomic Medicine"

# New "25 Students" class-size label (used by D15 and D16)
$ws.Range("D15").Value = "25 Students"

# New row 16 content
$ws.Range("A16").Value = "Cell Atlases"
$ws.Range("B16").Value = "MSc in Molecular Medicine, Genomic Medicine"

# New row 15 topic
$ws.Range("A15").Value = "RNA-sequencing (3 Lectures)"

# Update row 14's topic text to note it is now a full week course
$ws.Range("A14").Value = "Introduction to bioinformatics, transcriptomics and single cell RNA-sequencing (Full week)"

# New row 17 content (Trinity Single Cell Analysis Workshop)
$ws.Range("A17").Value = "Trinity Single Cell Analysis Workshop (Monthy Meeting)"
$ws.Range("B17").Value = "Open to Postgraduates and Post-Docs"
$ws.Range("D17").Value = "12 Members"

# D16 shares the same "25 Students" label as D15
$ws.Range("D16").Value = "25 Students"

# Year values for the three new rows
$ws.Range("C15").Value = 2025
$ws.Range("C16").Value = 2025
$ws.Range("C17").Value = 2025

# Match formatting (Arial font, style index 1) used by the rest of column A by copying
# format only (not value) from an existing column-A cell that already has that style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the view: scroll back to the top-left of the sheet and move the active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E18").Select() | Out-Null
